$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Collapse the "Egresos" detail block: remove the two extra detail rows
#    (old "Egreso 3"/"Egreso 4" helper rows) and the old "Total Egresos" row.
#    This shifts the old "Acumulado" row (15) up to row 12, keeping its
#    existing total-row formatting (fill + border) intact.
# ---------------------------------------------------------------------------
$ws.Range("A12:A14").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2) Add the new "Mes 4" column (column E), matching the width of the other
#    month columns (9.17 "characters" round-trips to the same stored
#    width="10" used by columns B-D).
# ---------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 9.17
$ws.Range("E1").Value = "Mes 4"

# ---------------------------------------------------------------------------
# 3) Rename the two remaining "Egreso" detail rows to "Ingenieros".
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "Ingenieros"
$ws.Range("A10").Value = "Ingenieros"

# ---------------------------------------------------------------------------
# 4) Row 11 used to be the "Egreso 3" detail row; turn it into the new
#    "Total Egresos" row, copying the total-row formatting from the
#    "Total Ingresos" row (row 7) so the fill/border matches.
# ---------------------------------------------------------------------------
$ws.Range("B7:D7").Copy()
$ws.Range("B11:D11").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A11").Value = "Total Egresos"

# ---------------------------------------------------------------------------
# 5) Update the monthly values across the sheet.
# ---------------------------------------------------------------------------
# -- Ingresos(*) detail --
$ws.Range("B3").Value = 12500
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 4488000
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

# -- Total Ingresos --
$ws.Range("B7").Value = 4500500
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

# -- Egresos(*) detail (now "Ingenieros" x2) --
$ws.Range("B9").Value = 8976
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 13464
$ws.Range("D10").Value = 0

# -- Total Egresos --
$ws.Range("B11").Value = 8976
$ws.Range("C11").Value = 13464
$ws.Range("D11").Value = 0

# -- Acumulado (old row 15, now shifted to row 12) --
$ws.Range("B12").Value = 4491524
$ws.Range("C12").Value = -13464
$ws.Range("D12").Value = 0

# ---------------------------------------------------------------------------
# 6) Fill in the new "Mes 4" column (all zeros, matching each row's styling).
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("E10").Value = 0

# Total-row cells in the new column need the same formatting as the rest of
# the corresponding total row.
$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("E7").Value = 0

$ws.Range("D11").Copy()
$ws.Range("E11").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("E11").Value = 0

$ws.Range("D12").Copy()
$ws.Range("E12").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("E12").Value = 0
